$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57, shifting existing rows 57..167 down to 58..168
$ws.Rows.Item(57).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Populate the newly inserted row 57 with its data
$ws.Cells.Item(57, 1).Value = 3
$ws.Cells.Item(57, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(57, 3).Value = "Coquimbo"
$ws.Cells.Item(57, 4).Value = 44645
$ws.Cells.Item(57, 5).Value = 5
$ws.Cells.Item(57, 6).Value = 100112030
$ws.Cells.Item(57, 7).Value = "Poroto granado"
$ws.Cells.Item(57, 8).Value = "Sin especificar"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 38
$ws.Cells.Item(57, 11).Value = 22000
$ws.Cells.Item(57, 12).Value = 22000
$ws.Cells.Item(57, 13).Value = 22000
$ws.Cells.Item(57, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(57, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(57, 16).Value = 880
$ws.Cells.Item(57, 17).Value = 25
$ws.Cells.Item(57, 18).Value = "Hortaliza"
